# Applies the refreshed cryptocurrency price/volume figures captured by the
# scheduled GitHub Actions run. Coin name/link/price/volume cells are updated
# in place; rows 47/48 additionally swap their ranking order (ApeXProtocol now
# outranks Fetch.AI).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '72.189.74'
$ws.Range("E2").Value = '  +4.18%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '4.035.53'
$ws.Range("E3").Value = '  +3.59%  '

# Row 4: TetherUSD
$ws.Range("E4").Value = '  -0.11%  '

# Row 5: BNB
$ws.Range("D5").Value = '''521.18'
$ws.Range("E5").Value = '  -1.37%  '

# Row 6: Solana
$ws.Range("D6").Value = '''147.22'
$ws.Range("E6").Value = '  +1.79%  '

# Row 7: XRP
$ws.Range("D7").Value = '''0.714'
$ws.Range("E7").Value = '  +16.34%  '

# Row 8: LidoStakedEther
$ws.Range("D8").Value = '4.027.27'
$ws.Range("E8").Value = '  +3.63%  '

# Row 9: USDC
$ws.Range("E9").Value = '  +0.14%  '

# Row 10: Cardano
$ws.Range("D10").Value = '''0.771'
$ws.Range("E10").Value = '  +7.31%  '

# Row 11: Dogecoin
$ws.Range("D11").Value = '''0.178'
$ws.Range("E11").Value = '  +3.85%  '

# Row 12: ShibaInu
$ws.Range("D12").Value = '''0.0000330'
$ws.Range("E12").Value = '  -0.80%  '

# Row 13: Avalanche
$ws.Range("D13").Value = '''48.57'
$ws.Range("E13").Value = '  +15.65%  '

# Row 14: Polkadot
$ws.Range("D14").Value = '''11.16'
$ws.Range("E14").Value = '  +8.94%  '

# Row 15: WrappedliquidstakedEther2.0
$ws.Range("D15").Value = '4.681.02'
$ws.Range("E15").Value = '  +3.48%  '

# Row 16: WrappedEther
$ws.Range("D16").Value = '4.069.32'
$ws.Range("E16").Value = '  +4.40%  '

# Row 17: Chainlink
$ws.Range("D17").Value = '''21.28'
$ws.Range("E17").Value = '  +7.74%  '

# Row 18: Uniswap
$ws.Range("D18").Value = '''14.26'
$ws.Range("E18").Value = '  +1.96%  '

# Row 19: Polygon
$ws.Range("E19").Value = '  +0.99%  '

# Row 20: TRON
$ws.Range("E20").Value = '  -0.42%  '

# Row 21: WrappedBTC
$ws.Range("D21").Value = '72.162.71'
$ws.Range("E21").Value = '  +4.16%  '

# Row 22: BitcoinCash
$ws.Range("D22").Value = '''445.18'
$ws.Range("E22").Value = '  +5.10%  '

# Row 23: Litecoin
$ws.Range("D23").Value = '''105.12'
$ws.Range("E23").Value = '  +19.85%  '

# Row 24: ImmutableX
$ws.Range("D24").Value = '''3.60'
$ws.Range("E24").Value = '  +6.13%  '

# Row 25: InternetComputer(DFINITY)
$ws.Range("D25").Value = '''15.30'
$ws.Range("E25").Value = '  +7.98%  '

# Row 26: PancakeSwap
$ws.Range("E26").Value = '  +0.91%  '

# Row 27: RenderToken
$ws.Range("D27").Value = '''11.50'
$ws.Range("E27").Value = '  +1.18%  '

# Row 28: Filecoin
$ws.Range("D28").Value = '''11.05'
$ws.Range("E28").Value = '  +4.25%  '

# Row 29: EthereumClassic
$ws.Range("D29").Value = '''37.89'
$ws.Range("E29").Value = '  +4.32%  '

# Row 30: LEO
$ws.Range("E30").Value = '  +2.36%  '

# Row 31: Toncoin
$ws.Range("D31").Value = '''3.31'
$ws.Range("E31").Value = '  +16.00%  '

# Row 32: Cosmos
$ws.Range("D32").Value = '''13.80'
$ws.Range("E32").Value = '  +4.53%  '

# Row 33: Hedera
$ws.Range("E33").Value = '  +3.40%  '

# Row 34: Bittensor
$ws.Range("D34").Value = '''677.64'
$ws.Range("E34").Value = '  -1.79%  '

# Row 35: NEARProtocol
$ws.Range("D35").Value = '''6.76'
$ws.Range("E35").Value = '  +14.46%  '

# Row 36: OKB
$ws.Range("D36").Value = '''67.94'
$ws.Range("E36").Value = '  -0.55%  '

# Row 37: InjectiveProtocol
$ws.Range("D37").Value = '''42.39'
$ws.Range("E37").Value = '  +6.18%  '

# Row 38: PEPE
$ws.Range("D38").Value = '0.0₃0864'
$ws.Range("E38").Value = '  +1.03%  '

# Row 39: TheGraph
$ws.Range("D39").Value = '''0.430'
$ws.Range("E39").Value = '  +0.38%  '

# Row 40: ThetaToken
$ws.Range("E40").Value = '  +6.61%  '

# Row 41: Kaspa
$ws.Range("D41").Value = '''0.153'
$ws.Range("E41").Value = '  +3.02%  '

# Row 42: Dai
$ws.Range("E42").Value = '  +0.22%  '

# Row 43: VeChain
$ws.Range("D43").Value = '''0.0502'
$ws.Range("E43").Value = '  +4.15%  '

# Row 45: WEMIXToken
$ws.Range("D45").Value = '''3.21'
$ws.Range("E45").Value = '  -0.75%  '

# Row 46: Stellar
$ws.Range("D46").Value = '''0.157'
$ws.Range("E46").Value = '  +12.00%  '

# Row 47: Fetch.AI
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").Value = '''3.54'
$ws.Range("E47").Value = '  +3.71%  '

# Row 48: ApeXProtocol
$ws.Range("B48").Value = 'Fetch.AI'
$ws.Range("C48").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D48").Value = '''2.71'
$ws.Range("E48").Value = '  -2.52%  '

# Row 49: THORChain
$ws.Range("D49").Value = '''9.56'
$ws.Range("E49").Value = '  +11.54%  '

# Row 50: Stacks
$ws.Range("D50").Value = '''3.08'
$ws.Range("E50").Value = '  +2.35%  '

# Row 51: LidoDAOToken
$ws.Range("E51").Value = '  +2.97%  '
